$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11, column B currently holds the text "R40" (row label for the
# fourth rule). Update it to the new label "1", stored as text so the
# shared-string table gets a new string entry (matching the target
# workbook which keeps t="s" on B11).
$ws.Range("B11").Value = "1"
